# Scale the Assets sheet parse output from 10 objects to 20 objects.
# Appends 36 new "image" rows (rows 44-79) for snowboard objects 10-19,
# following the same variant/image layout as the existing rows, and
# updates the sheet dimension + ignoredErrors ranges to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("image","snowboard1010","SNW-2015-01","SNW-10-01","SNW-Image10-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-10-01/0.jpg","700","700"),
    @("image","snowboard1010","SNW-2015-01","SNW-10-01","SNW-Image10-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-10-01/1.jpg","700","700"),
    @("image","snowboard1010","SNW-2015-01","SNW-10-01","SNW-Image10-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-10-01/2.jpg","700","700"),
    @("image","snowboard1010","SNW-2015-01","SNW-10-01","SNW-Image10-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-10-01/3.jpg","700","700"),
    @("image","snowboard1010","SNW-2015-01","SNW-10-01","SNW-Image10-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-10-01/4.jpg","700","700"),
    @("image","snowboard1011","SNW-3475-01","SNW-11-01","SNW-Image11-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-11-01/0.jpg","700","700"),
    @("image","snowboard1012","SNW-2784-01","SNW-12-01","SNW-Image12-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-12-01/0.jpg","700","700"),
    @("image","snowboard1012","SNW-2784-01","SNW-12-01","SNW-Image12-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-12-01/1.jpg","700","700"),
    @("image","snowboard1012","SNW-2784-01","SNW-12-01","SNW-Image12-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-12-01/2.jpg","700","700"),
    @("image","snowboard1012","SNW-2784-01","SNW-12-01","SNW-Image12-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-12-01/3.jpg","700","700"),
    @("image","snowboard1012","SNW-2784-01","SNW-12-01","SNW-Image12-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-12-01/4.jpg","700","700"),
    @("image","snowboard1013","SNW-3275-01","SNW-13-01","SNW-Image13-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-13-01/0.jpg","700","700"),
    @("image","snowboard1013","SNW-3275-01","SNW-13-01","SNW-Image13-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-13-01/1.jpg","700","700"),
    @("image","snowboard1013","SNW-3275-01","SNW-13-01","SNW-Image13-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-13-01/2.jpg","700","700"),
    @("image","snowboard1013","SNW-3275-01","SNW-13-01","SNW-Image13-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-13-01/3.jpg","700","700"),
    @("image","snowboard1013","SNW-3275-01","SNW-13-01","SNW-Image13-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-13-01/4.jpg","700","700"),
    @("image","snowboard1014","SNW-3655-01","SNW-14-01","SNW-Image14-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-14-01/0.jpg","700","700"),
    @("image","snowboard1015","SNW-3820-01","SNW-15-01","SNW-Image15-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-15-01/0.jpg","700","700"),
    @("image","snowboard1016","SNW-2578-01","SNW-16-01","SNW-Image16-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-16-01/0.jpg","700","700"),
    @("image","snowboard1016","SNW-2578-01","SNW-16-01","SNW-Image16-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-16-01/1.jpg","700","700"),
    @("image","snowboard1016","SNW-2578-01","SNW-16-01","SNW-Image16-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-16-01/2.jpg","700","700"),
    @("image","snowboard1016","SNW-2578-01","SNW-16-01","SNW-Image16-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-16-01/3.jpg","700","700"),
    @("image","snowboard1016","SNW-2578-01","SNW-16-01","SNW-Image16-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-16-01/4.jpg","700","700"),
    @("image","snowboard1017","SNW-3321-01","SNW-17-01","SNW-Image17-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-17-01/0.jpg","700","700"),
    @("image","snowboard1017","SNW-3321-01","SNW-17-01","SNW-Image17-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-17-01/1.jpg","700","700"),
    @("image","snowboard1017","SNW-3321-01","SNW-17-01","SNW-Image17-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-17-01/2.jpg","700","700"),
    @("image","snowboard1017","SNW-3321-01","SNW-17-01","SNW-Image17-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-17-01/3.jpg","700","700"),
    @("image","snowboard1017","SNW-3321-01","SNW-17-01","SNW-Image17-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-17-01/4.jpg","700","700"),
    @("image","snowboard1018","SNW-3658-01","SNW-18-01","SNW-Image18-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-18-01/0.jpg","700","700"),
    @("image","snowboard1018","SNW-3658-01","SNW-18-01","SNW-Image18-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-18-01/1.jpg","700","700"),
    @("image","snowboard1018","SNW-3658-01","SNW-18-01","SNW-Image18-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-18-01/2.jpg","700","700"),
    @("image","snowboard1018","SNW-3658-01","SNW-18-01","SNW-Image18-04","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-18-01/3.jpg","700","700"),
    @("image","snowboard1018","SNW-3658-01","SNW-18-01","SNW-Image18-05","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-18-01/4.jpg","700","700"),
    @("image","snowboard1019","SNW-2798-01","SNW-19-01","SNW-Image19-01","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-19-01/0.jpg","700","700"),
    @("image","snowboard1019","SNW-2798-01","SNW-19-01","SNW-Image19-02","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-19-01/1.jpg","700","700"),
    @("image","snowboard1019","SNW-2798-01","SNW-19-01","SNW-Image19-03","https://raw.githubusercontent.com/ksarise/parser/main/assets/SNW-19-01/2.jpg","700","700")
)

$startRow = 44
$r = $startRow
foreach ($row in $newRows) {
    for ($c = 1; $c -le 8; $c++) {
        $val = $row[$c - 1]
        if ($c -ge 7) {
            # Columns G/H hold numeric-looking values ("700") that must be
            # stored as text (matching the source data's string typing), so
            # prefix with an apostrophe to force text entry instead of
            # Excel auto-coercing them into numbers.
            $ws.Cells.Item($r, $c).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
    $r++
}

$endRow = $r - 1

$ws.Range("A1:H$endRow").Select()
